# Project Evaluation Data.xlsx - add analysis (IPD benchmark row, accuracy /
# difference columns, correlation summary, evaluation & analysis notes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 1 header row: re-label M1, add headers for the new N:T columns.
# ---------------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true

$ws.Range("M1").Value = "% Change 4/3-4/10"
$ws.Range("M1").Font.Bold = $true

$ws.Range("N1").Value = "% positivity (all tweets)"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").Font.Color = 0
$ws.Range("N1").Font.Name = "Calibri"

$ws.Range("O1").Value = "Accurate?"
$ws.Range("O1").Font.Bold = $true

$ws.Range("P1").Value = "Difference"
$ws.Range("P1").Font.Bold = $true

$ws.Range("Q1").Value = "% positivity ( > 1000)"
$ws.Range("Q1").Font.Bold = $true
$ws.Range("Q1").Font.Color = 0
$ws.Range("Q1").Font.Name = "Calibri"

$ws.Range("R1").Value = "Accurate?"
$ws.Range("R1").Font.Bold = $true

$ws.Range("S1").Value = "% positivity ( > 10000)"
$ws.Range("S1").Font.Bold = $true
$ws.Range("S1").Font.Color = 0
$ws.Range("S1").Font.Name = "Calibri"

$ws.Range("T1").Value = "Accurate?"
$ws.Range("T1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Rows 2-6 (GM, Ford, VW, Toyota, Tesla): weekly % change formula in M,
#    plus tweet-positivity inputs + accuracy/difference columns N:U.
# ---------------------------------------------------------------------------

# -- Row 2 (GM) --
$ws.Range("M2").Formula = "=(L2-B2)/B2"
$ws.Range("M2").NumberFormat = "0.000%"
$ws.Range("N2").Value = -0.0163
$ws.Range("N2").NumberFormat = "0.000%"
$ws.Range("O2").Formula = '=IF(OR(AND(N2>0,$M2>0),AND(N2<0,$M2<0)),"YES","NO")'
$ws.Range("P2").Formula = "=ABS(N2-M2)"
$ws.Range("P2").NumberFormat = "0.000%"
$ws.Range("Q2").Value = -0.00331
$ws.Range("Q2").NumberFormat = "0.000%"
$ws.Range("R2").Formula = '=IF(OR(AND(Q2>0,$M2>0),AND(Q2<0,$M2<0)),"YES","NO")'
$ws.Range("S2").Value = 0.0146
$ws.Range("S2").NumberFormat = "0.000%"
$ws.Range("T2").Formula = '=IF(OR(AND(S2>0,$M2>0),AND(S2<0,$M2<0)),"YES","NO")'
$ws.Range("U2").NumberFormat = "0.000%"

# -- Row 3 (Ford) --
$ws.Range("M3").Formula = "=(L3-B3)/B3"
$ws.Range("M3").NumberFormat = "0.000%"
$ws.Range("N3").Value = 0.0075
$ws.Range("N3").NumberFormat = "0.000%"
$ws.Range("O3").Formula = '=IF(OR(AND(N3>0,$M3>0),AND(N3<0,$M3<0)),"YES","NO")'
$ws.Range("P3").Formula = "=ABS(N3-M3)"
$ws.Range("P3").NumberFormat = "0.000%"
$ws.Range("Q3").Value = 0.00635
$ws.Range("Q3").NumberFormat = "0.000%"
$ws.Range("R3").Formula = '=IF(OR(AND(Q3>0,$M3>0),AND(Q3<0,$M3<0)),"YES","NO")'
$ws.Range("S3").Value = -0.01663
$ws.Range("S3").NumberFormat = "0.000%"
$ws.Range("T3").Formula = '=IF(OR(AND(S3>0,$M3>0),AND(S3<0,$M3<0)),"YES","NO")'
$ws.Range("U3").NumberFormat = "0.000%"

# -- Row 4 (VW) --
$ws.Range("M4").Formula = "=(L4-B4)/B4"
$ws.Range("M4").NumberFormat = "0.000%"
$ws.Range("N4").Value = 0.0088
$ws.Range("N4").NumberFormat = "0.000%"
$ws.Range("O4").Formula = '=IF(OR(AND(N4>0,$M4>0),AND(N4<0,$M4<0)),"YES","NO")'
$ws.Range("P4").Formula = "=ABS(N4-M4)"
$ws.Range("P4").NumberFormat = "0.000%"
$ws.Range("Q4").Value = 0.02356
$ws.Range("Q4").NumberFormat = "0.000%"
$ws.Range("R4").Formula = '=IF(OR(AND(Q4>0,$M4>0),AND(Q4<0,$M4<0)),"YES","NO")'
$ws.Range("S4").Value = 0.0349
$ws.Range("S4").NumberFormat = "0.000%"
$ws.Range("T4").Formula = '=IF(OR(AND(S4>0,$M4>0),AND(S4<0,$M4<0)),"YES","NO")'
$ws.Range("U4").NumberFormat = "0.000%"

# -- Row 5 (Toyota) --
$ws.Range("M5").Formula = "=(L5-B5)/B5"
$ws.Range("M5").NumberFormat = "0.000%"
$ws.Range("N5").Value = 0.0072
$ws.Range("N5").NumberFormat = "0.000%"
$ws.Range("O5").Formula = '=IF(OR(AND(N5>0,$M5>0),AND(N5<0,$M5<0)),"YES","NO")'
$ws.Range("P5").Formula = "=ABS(N5-M5)"
$ws.Range("P5").NumberFormat = "0.000%"
$ws.Range("Q5").Value = 0.00512
$ws.Range("Q5").NumberFormat = "0.000%"
$ws.Range("R5").Formula = '=IF(OR(AND(Q5>0,$M5>0),AND(Q5<0,$M5<0)),"YES","NO")'
$ws.Range("S5").Value = 0.04475
$ws.Range("S5").NumberFormat = "0.000%"
$ws.Range("T5").Formula = '=IF(OR(AND(S5>0,$M5>0),AND(S5<0,$M5<0)),"YES","NO")'
$ws.Range("U5").NumberFormat = "0.000%"

# -- Row 6 (Tesla) --
$ws.Range("M6").Formula = "=(L6-B6)/B6"
$ws.Range("M6").NumberFormat = "0.000%"
$ws.Range("N6").Value = 0.00654
$ws.Range("N6").NumberFormat = "0.000%"
$ws.Range("O6").Formula = '=IF(OR(AND(N6>0,$M6>0),AND(N6<0,$M6<0)),"YES","NO")'
$ws.Range("P6").Formula = "=ABS(N6-M6)"
$ws.Range("P6").NumberFormat = "0.000%"
$ws.Range("Q6").Value = 0.01231
$ws.Range("Q6").NumberFormat = "0.000%"
$ws.Range("R6").Formula = '=IF(OR(AND(Q6>0,$M6>0),AND(Q6<0,$M6<0)),"YES","NO")'
$ws.Range("S6").Value = 0.02365
$ws.Range("S6").NumberFormat = "0.000%"
$ws.Range("T6").Formula = '=IF(OR(AND(S6>0,$M6>0),AND(S6<0,$M6<0)),"YES","NO")'
$ws.Range("U6").NumberFormat = "0.000%"

# ---------------------------------------------------------------------------
# 3. Row 7: averages across the 5 stocks.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Average of 5 stocks:"
$ws.Range("M7").Formula = "=AVERAGE(M2:M6)"
$ws.Range("M7").NumberFormat = "0.000%"
$ws.Range("N7").Formula = "=AVERAGE(N2:N6)"
$ws.Range("N7").NumberFormat = "0.000%"
$ws.Range("O7").NumberFormat = "0.000%"
$ws.Range("P7").NumberFormat = "0.000%"
$ws.Range("Q7").Formula = "=AVERAGE(Q2:Q6)"
$ws.Range("Q7").NumberFormat = "0.000%"
$ws.Range("R7").NumberFormat = "0.000%"
$ws.Range("S7").Formula = "=AVERAGE(S2:S6)"
$ws.Range("S7").NumberFormat = "0.000%"

# ---------------------------------------------------------------------------
# 4. Row 9: IPD (Auto Industry ETF) benchmark data, same layout as rows 2-6.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "IPD (Auto Industry ETF)"
$ws.Range("B9").Value = 37.96
$ws.Range("C9").Formula = "=(D9-B9)/B9*100"
$ws.Range("D9").Value = 37.87
$ws.Range("E9").Formula = "=(F9-D9)/D9*100"
$ws.Range("F9").Value = 37.7
$ws.Range("G9").Formula = "=(H9-F9)/F9*100"
$ws.Range("H9").Value = 37.65
$ws.Range("I9").Formula = "=(J9-H9)/H9*100"
$ws.Range("J9").Value = 37.59
$ws.Range("K9").Formula = "=(L9-J9)/J9*100"
$ws.Range("L9").Value = 37.59
$ws.Range("M9").Formula = "=(L9-B9)/B9"
$ws.Range("M9").NumberFormat = "0.000%"

# ---------------------------------------------------------------------------
# 5. Row 11: correlation of weekly % change against each positivity metric.
# ---------------------------------------------------------------------------
$ws.Range("M11").Value = "Correlation:"
$ws.Range("N11").Formula = "=CORREL(M2:M6,N2:N6)"
$ws.Range("Q11").Formula = "=CORREL(M2:M6,Q2:Q6)"
$ws.Range("S11").Formula = "=CORREL(M2:M6,S2:S6)"

# ---------------------------------------------------------------------------
# 6. Rows 16-18: Evaluation notes.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Evaluation:"
$ws.Range("A17").Value = "'-Correlation between % change (week) and % positivity"
$ws.Range("A18").Value = "'-Accuracy (whether sign of our prediction matches weekly price movement)"

# ---------------------------------------------------------------------------
# 7. Rows 22-25: Analysis notes.
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "Analysis:"
$ws.Range("A23").Value = "'-Follower count(s) correlations"
$ws.Range("A24").Value = "'-Industry correlation (auto ETF) with each stock"
$ws.Range("A25").Value = '''-"The greater the % change, the further our model is"'

# ---------------------------------------------------------------------------
# 8. Column widths matching the authored layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.1666666667
$ws.Columns.Item(13).ColumnWidth = 15.8307291667
$ws.Columns.Item(14).ColumnWidth = 19.6666666667
$ws.Columns.Item(17).ColumnWidth = 17.4986979167
$ws.Columns.Item(18).ColumnWidth = 17.4986979167
$ws.Columns.Item(19).ColumnWidth = 17.4986979167
$ws.Columns.Item(21).ColumnWidth = 11.3307291667

# ---------------------------------------------------------------------------
# 9. Selection mirrors where the author left off.
# ---------------------------------------------------------------------------
$ws.Range("A26").Select()
